$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32 (sheet ALC)
$ws.Cells.Item(32, 8).Value = 2583.5715  # H32: was 2474.6667
$ws.Cells.Item(32, 9).Value = 1758.8  # I32: was 1623.1666
$ws.Cells.Item(32, 10).Value = 3041.7778  # J32: was 3042.3333
$ws.Cells.Item(32, 11).Value = 1758.8  # K32: was 1623.1666
$ws.Cells.Item(32, 12).Value = 3041.7778  # L32: was 3042.3333
$ws.Cells.Item(32, 13).Value = -1432.8  # M32: was -1297.1666
$ws.Cells.Item(32, 14).Value = -3693.7778  # N32: was -3694.3333

# Row 100 (sheet ALC)
$ws.Cells.Item(100, 8).Value = 1270.5  # H100: was 1420
$ws.Cells.Item(100, 9).Value = 1174.5  # I100: was 0
$ws.Cells.Item(100, 10).Value = 1366.5  # J100: was 1420
$ws.Cells.Item(100, 11).Value = 1174.5  # K100: was 0
$ws.Cells.Item(100, 12).Value = 1366.5  # L100: was 1420
$ws.Cells.Item(100, 13).Value = -633.5  # M100: was None
$ws.Cells.Item(100, 14).Value = -2448.5  # N100: was -2502

# Row 104 (sheet ALC)
$ws.Cells.Item(104, 8).Value = 148.8  # H104: was 140.66667
$ws.Cells.Item(104, 9).Value = 148.8  # I104: was 140.66667
$ws.Cells.Item(104, 11).Value = 446.4  # K104: was 422.00001
$ws.Cells.Item(104, 13).Value = 1300.6  # M104: was 1324.99999

# Row 105 (sheet ALC)
$ws.Cells.Item(105, 8).Value = 32499.5  # H105: was 0
$ws.Cells.Item(105, 10).Value = 32499.5  # J105: was 0
$ws.Cells.Item(105, 12).Value = 32499.5  # L105: was 0
$ws.Cells.Item(105, 14).Value = -39487.5  # N105: was None

# Row 106 (sheet ALC)
$ws.Cells.Item(106, 8).Value = 41177.6  # H106: was 41197.6
$ws.Cells.Item(106, 9).Value = 45222  # I106: was 45247
$ws.Cells.Item(106, 11).Value = 45222  # K106: was 45247
$ws.Cells.Item(106, 13).Value = -44591  # M106: was -44616

# Row 115 (sheet ALC)
$ws.Cells.Item(115, 8).Value = 184.33333  # H115: was 394.66666
$ws.Cells.Item(115, 9).Value = 184.33333  # I115: was 394.66666
$ws.Cells.Item(115, 11).Value = 552.99999  # K115: was 1183.99998
$ws.Cells.Item(115, 13).Value = 1014.00001  # M115: was 383.0000199999999

# Row 132 (sheet ALC)
$ws.Cells.Item(132, 8).Value = 1193.5  # H132: was 1288.65
$ws.Cells.Item(132, 9).Value = 1193.5  # I132: was 1288.65
$ws.Cells.Item(132, 11).Value = 3580.5  # K132: was 3865.95
$ws.Cells.Item(132, 13).Value = -1050.5  # M132: was -1335.95

# Row 137 (sheet ALC)
$ws.Cells.Item(137, 8).Value = 1298.8235  # H137: was 1463.875
$ws.Cells.Item(137, 9).Value = 1335.9166  # I137: was 1347.9166
$ws.Cells.Item(137, 10).Value = 1209.8  # J137: was 1811.75
$ws.Cells.Item(137, 11).Value = 4007.7498  # K137: was 4043.7498
$ws.Cells.Item(137, 12).Value = 3629.4  # L137: was 5435.25
$ws.Cells.Item(137, 13).Value = -1457.7498  # M137: was -1493.7498
$ws.Cells.Item(137, 14).Value = -8729.4  # N137: was -10535.25

$ws = $wb.Worksheets.Item("ARM")
# Row 4 (sheet ARM)
$ws.Cells.Item(4, 8).Value = 2119.6  # H4: was 3166.3333
$ws.Cells.Item(4, 9).Value = 599  # I4: was 0
$ws.Cells.Item(4, 10).Value = 2499.75  # J4: was 3166.3333
$ws.Cells.Item(4, 11).Value = 599  # K4: was 0
$ws.Cells.Item(4, 12).Value = 2499.75  # L4: was 3166.3333
$ws.Cells.Item(4, 13).Value = -483  # M4: was None
$ws.Cells.Item(4, 14).Value = -2731.75  # N4: was -3398.3333

# Row 5 (sheet ARM)
$ws.Cells.Item(5, 8).Value = 238.6  # H5: was 260.55554
$ws.Cells.Item(5, 9).Value = 318.33334  # I5: was 373.8
$ws.Cells.Item(5, 11).Value = 318.33334  # K5: was 373.8
$ws.Cells.Item(5, 13).Value = -206.33334  # M5: was -261.8

# Row 32 (sheet ARM)
$ws.Cells.Item(32, 8).Value = 2991.4443  # H32: was 3000.0278
$ws.Cells.Item(32, 9).Value = 3089.9333  # I32: was 3100.2334
$ws.Cells.Item(32, 11).Value = 3089.9333  # K32: was 3100.2334
$ws.Cells.Item(32, 13).Value = -2802.9333  # M32: was -2813.2334

# Row 132 (sheet ARM)
$ws.Cells.Item(132, 8).Value = 3254.6155  # H132: was 2360.5264
$ws.Cells.Item(132, 9).Value = 3191.3333  # I132: was 2268.6667
$ws.Cells.Item(132, 11).Value = 9573.999899999999  # K132: was 6806.000100000001
$ws.Cells.Item(132, 13).Value = -7043.999899999999  # M132: was -4276.000100000001

$ws = $wb.Worksheets.Item("BSM")
# Row 4 (sheet BSM)
$ws.Cells.Item(4, 8).Value = 238.6  # H4: was 260.55554
$ws.Cells.Item(4, 9).Value = 318.33334  # I4: was 373.8
$ws.Cells.Item(4, 11).Value = 318.33334  # K4: was 373.8
$ws.Cells.Item(4, 13).Value = -203.33334  # M4: was -258.8

# Row 37 (sheet BSM)
$ws.Cells.Item(37, 8).Value = 1090.6666  # H37: was 848.25
$ws.Cells.Item(37, 9).Value = 1090.6666  # I37: was 848.25
$ws.Cells.Item(37, 11).Value = 1090.6666  # K37: was 848.25
$ws.Cells.Item(37, 13).Value = -953.6666  # M37: was -711.25

# Row 94 (sheet BSM)
$ws.Cells.Item(94, 8).Value = 1065.6666  # H94: was 1078.1
$ws.Cells.Item(94, 9).Value = 663.5  # I94: was 799
$ws.Cells.Item(94, 10).Value = 1387.4  # J94: was 1496.75
$ws.Cells.Item(94, 11).Value = 663.5  # K94: was 799
$ws.Cells.Item(94, 12).Value = 1387.4  # L94: was 1496.75
$ws.Cells.Item(94, 13).Value = -212.5  # M94: was -348
$ws.Cells.Item(94, 14).Value = -2289.4  # N94: was -2398.75

# Row 105 (sheet BSM)
$ws.Cells.Item(105, 8).Value = 2609.2666  # H105: was 2807.6155
$ws.Cells.Item(105, 9).Value = 2494.5454  # I105: was 2755.5557
$ws.Cells.Item(105, 11).Value = 2494.5454  # K105: was 2755.5557
$ws.Cells.Item(105, 13).Value = -747.5454  # M105: was -1008.5557

# Row 107 (sheet BSM)
$ws.Cells.Item(107, 8).Value = 2341.0967  # H107: was 2450.9312
$ws.Cells.Item(107, 9).Value = 1978.12  # I107: was 2085.0435
$ws.Cells.Item(107, 11).Value = 1978.12  # K107: was 2085.0435
$ws.Cells.Item(107, 13).Value = -58.11999999999989  # M107: was -165.0435000000002

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (sheet CRP)
$ws.Cells.Item(16, 8).Value = 2664.889  # H16: was 2344.818
$ws.Cells.Item(16, 9).Value = 2997.8  # I16: was 2399.7144
$ws.Cells.Item(16, 11).Value = 2997.8  # K16: was 2399.7144
$ws.Cells.Item(16, 13).Value = -2710.8  # M16: was -2112.7144

# Row 58 (sheet CRP)
$ws.Cells.Item(58, 8).Value = 1350  # H58: was 1525.3889
$ws.Cells.Item(58, 9).Value = 1123.1578  # I58: was 1289.6154
$ws.Cells.Item(58, 10).Value = 2427.5  # J58: was 2138.4
$ws.Cells.Item(58, 11).Value = 1123.1578  # K58: was 1289.6154
$ws.Cells.Item(58, 12).Value = 2427.5  # L58: was 2138.4
$ws.Cells.Item(58, 13).Value = -920.1578  # M58: was -1086.6154
$ws.Cells.Item(58, 14).Value = -2833.5  # N58: was -2544.4

# Row 113 (sheet CRP)
$ws.Cells.Item(113, 8).Value = 2664.889  # H113: was 2344.818
$ws.Cells.Item(113, 9).Value = 2997.8  # I113: was 2399.7144
$ws.Cells.Item(113, 11).Value = 2997.8  # K113: was 2399.7144
$ws.Cells.Item(113, 13).Value = -827.8000000000002  # M113: was -229.7143999999998

# Row 122 (sheet CRP)
$ws.Cells.Item(122, 8).Value = 2828.2144  # H122: was 3357.3
$ws.Cells.Item(122, 9).Value = 2264  # I122: was 2510.8572
$ws.Cells.Item(122, 10).Value = 3843.8  # J122: was 5332.3335
$ws.Cells.Item(122, 11).Value = 6792  # K122: was 7532.571599999999
$ws.Cells.Item(122, 12).Value = 11531.4  # L122: was 15997.0005
$ws.Cells.Item(122, 13).Value = -4342  # M122: was -5082.571599999999
$ws.Cells.Item(122, 14).Value = -16431.4  # N122: was -20897.0005

# Row 136 (sheet CRP)
$ws.Cells.Item(136, 8).Value = 1350  # H136: was 1525.3889
$ws.Cells.Item(136, 9).Value = 1123.1578  # I136: was 1289.6154
$ws.Cells.Item(136, 10).Value = 2427.5  # J136: was 2138.4
$ws.Cells.Item(136, 11).Value = 3369.4734  # K136: was 3868.8462
$ws.Cells.Item(136, 12).Value = 7282.5  # L136: was 6415.200000000001
$ws.Cells.Item(136, 13).Value = -819.4733999999999  # M136: was -1318.8462
$ws.Cells.Item(136, 14).Value = -12382.5  # N136: was -11515.2

$ws = $wb.Worksheets.Item("CUL")
# Row 34 (sheet CUL)
$ws.Cells.Item(34, 8).Value = 631.6667  # H34: was 3465
$ws.Cells.Item(34, 10).Value = 1500  # J34: was 10000
$ws.Cells.Item(34, 12).Value = 4500  # L34: was 30000
$ws.Cells.Item(34, 14).Value = -4668  # N34: was -30168

# Row 69 (sheet CUL)
$ws.Cells.Item(69, 8).Value = 2394.3333  # H69: was 2454.9
$ws.Cells.Item(69, 10).Value = 2992.8572  # J69: was 2993.75
$ws.Cells.Item(69, 12).Value = 8978.571599999999  # L69: was 8981.25
$ws.Cells.Item(69, 14).Value = -10600.5716  # N69: was -10603.25

# Row 72 (sheet CUL)
$ws.Cells.Item(72, 8).Value = 2394.3333  # H72: was 2454.9
$ws.Cells.Item(72, 10).Value = 2992.8572  # J72: was 2993.75
$ws.Cells.Item(72, 12).Value = 26935.7148  # L72: was 26943.75
$ws.Cells.Item(72, 14).Value = -35047.7148  # N72: was -35055.75

# Row 107 (sheet CUL)
$ws.Cells.Item(107, 8).Value = 410  # H107: was 487.5
$ws.Cells.Item(107, 10).Value = 575  # J107: was 812.5
$ws.Cells.Item(107, 12).Value = 1725  # L107: was 2437.5
$ws.Cells.Item(107, 14).Value = -5565  # N107: was -6277.5

# Row 118 (sheet CUL)
$ws.Cells.Item(118, 8).Value = 2082.6086  # H118: was 1992.375
$ws.Cells.Item(118, 9).Value = 2082.6086  # I118: was 1992.375
$ws.Cells.Item(118, 11).Value = 6247.825800000001  # K118: was 5977.125
$ws.Cells.Item(118, 13).Value = -5004.825800000001  # M118: was -4734.125

# Row 139 (sheet CUL)
$ws.Cells.Item(139, 8).Value = 753.3333  # H139: was 721.875
$ws.Cells.Item(139, 9).Value = 630.5  # I139: was 721.875
$ws.Cells.Item(139, 10).Value = 999  # J139: was 0
$ws.Cells.Item(139, 11).Value = 1891.5  # K139: was 2165.625
$ws.Cells.Item(139, 12).Value = 2997  # L139: was 0
$ws.Cells.Item(139, 13).Value = 3248.5  # M139: was 2974.375
$ws.Cells.Item(139, 14).Value = -13277  # N139: was None

$ws = $wb.Worksheets.Item("GSM")
# Row 122 (sheet GSM)
$ws.Cells.Item(122, 8).Value = 3457.3684  # H122: was 3570.1765
$ws.Cells.Item(122, 9).Value = 3558.0833  # I122: was 3770
$ws.Cells.Item(122, 11).Value = 10674.2499  # K122: was 11310
$ws.Cells.Item(122, 13).Value = -8224.249899999999  # M122: was -8860

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (sheet LTW)
$ws.Cells.Item(7, 8).Value = 7872.45  # H7: was 7941.722
$ws.Cells.Item(7, 9).Value = 5692.4  # I7: was 4992.4
$ws.Cells.Item(7, 10).Value = 8599.134  # J7: was 9076.076999999999
$ws.Cells.Item(7, 11).Value = 5692.4  # K7: was 4992.4
$ws.Cells.Item(7, 12).Value = 8599.134  # L7: was 9076.076999999999
$ws.Cells.Item(7, 13).Value = -5580.4  # M7: was -4880.4
$ws.Cells.Item(7, 14).Value = -8823.134  # N7: was -9300.076999999999

# Row 22 (sheet LTW)
$ws.Cells.Item(22, 8).Value = 1974.8334  # H22: was 1841.6666
$ws.Cells.Item(22, 9).Value = 1949.6666  # I22: was 1683.3334
$ws.Cells.Item(22, 11).Value = 1949.6666  # K22: was 1683.3334
$ws.Cells.Item(22, 13).Value = -1654.6666  # M22: was -1388.3334

# Row 27 (sheet LTW)
$ws.Cells.Item(27, 8).Value = 1974.8334  # H27: was 1841.6666
$ws.Cells.Item(27, 9).Value = 1949.6666  # I27: was 1683.3334
$ws.Cells.Item(27, 11).Value = 1949.6666  # K27: was 1683.3334
$ws.Cells.Item(27, 13).Value = -1842.6666  # M27: was -1576.3334

# Row 40 (sheet LTW)
$ws.Cells.Item(40, 8).Value = 3682.3333  # H40: was 3377.7144
$ws.Cells.Item(40, 9).Value = 3254.5557  # I40: was 2944.6365
$ws.Cells.Item(40, 11).Value = 3254.5557  # K40: was 2944.6365
$ws.Cells.Item(40, 13).Value = -3118.5557  # M40: was -2808.6365

# Row 93 (sheet LTW)
$ws.Cells.Item(93, 8).Value = 3259.8  # H93: was 3800
$ws.Cells.Item(93, 9).Value = 3833  # I93: was 0
$ws.Cells.Item(93, 10).Value = 2400  # J93: was 3800
$ws.Cells.Item(93, 11).Value = 3833  # K93: was 0
$ws.Cells.Item(93, 12).Value = 2400  # L93: was 3800
$ws.Cells.Item(93, 13).Value = -2585  # M93: was None
$ws.Cells.Item(93, 14).Value = -4896  # N93: was -6296

# Row 122 (sheet LTW)
$ws.Cells.Item(122, 8).Value = 8078.737  # H122: was 7258.864
$ws.Cells.Item(122, 9).Value = 8566.6  # I122: was 7483.222
$ws.Cells.Item(122, 11).Value = 25699.8  # K122: was 22449.666
$ws.Cells.Item(122, 13).Value = -23249.8  # M122: was -19999.666

# Row 126 (sheet LTW)
$ws.Cells.Item(126, 8).Value = 7872.45  # H126: was 7941.722
$ws.Cells.Item(126, 9).Value = 5692.4  # I126: was 4992.4
$ws.Cells.Item(126, 10).Value = 8599.134  # J126: was 9076.076999999999
$ws.Cells.Item(126, 11).Value = 17077.2  # K126: was 14977.2
$ws.Cells.Item(126, 12).Value = 25797.402  # L126: was 27228.231
$ws.Cells.Item(126, 13).Value = -14607.2  # M126: was -12507.2
$ws.Cells.Item(126, 14).Value = -30737.402  # N126: was -32168.231

$ws = $wb.Worksheets.Item("WVR")
# Row 81 (sheet WVR)
$ws.Cells.Item(81, 8).Value = 2000790.2  # H81: was 1000609
$ws.Cells.Item(81, 9).Value = 987.25  # I81: was 648.5
$ws.Cells.Item(81, 10).Value = 10000002  # J81: was 5000451
$ws.Cells.Item(81, 11).Value = 1974.5  # K81: was 1297
$ws.Cells.Item(81, 12).Value = 20000004  # L81: was 10000902
$ws.Cells.Item(81, 13).Value = -913.5  # M81: was -236
$ws.Cells.Item(81, 14).Value = -20002126  # N81: was -10003024

# Row 82 (sheet WVR)
$ws.Cells.Item(82, 8).Value = 0  # H82: was 40301
$ws.Cells.Item(82, 10).Value = 0  # J82: was 40301
$ws.Cells.Item(82, 12).Value = 0  # L82: was 40301
$ws.Cells.Item(82, 14).ClearContents()  # N82: was -41067

# Row 84 (sheet WVR)
$ws.Cells.Item(84, 8).Value = 2000790.2  # H84: was 1000609
$ws.Cells.Item(84, 9).Value = 987.25  # I84: was 648.5
$ws.Cells.Item(84, 10).Value = 10000002  # J84: was 5000451
$ws.Cells.Item(84, 11).Value = 9872.5  # K84: was 6485
$ws.Cells.Item(84, 12).Value = 100000020  # L84: was 50004510
$ws.Cells.Item(84, 13).Value = -4568.5  # M84: was -1181
$ws.Cells.Item(84, 14).Value = -100010628  # N84: was -50015118

# Row 85 (sheet WVR)
$ws.Cells.Item(85, 8).Value = 0  # H85: was 40301
$ws.Cells.Item(85, 10).Value = 0  # J85: was 40301
$ws.Cells.Item(85, 12).Value = 0  # L85: was 40301
$ws.Cells.Item(85, 14).ClearContents()  # N85: was -42953

# Row 96 (sheet WVR)
$ws.Cells.Item(96, 8).Value = 997  # H96: was 999
$ws.Cells.Item(96, 9).Value = 997  # I96: was 999
$ws.Cells.Item(96, 11).Value = 997  # K96: was 999
$ws.Cells.Item(96, 13).Value = 376  # M96: was 374

# Row 98 (sheet WVR)
$ws.Cells.Item(98, 8).Value = 15999.5  # H98: was 15766.333
$ws.Cells.Item(98, 10).Value = 15999.5  # J98: was 15766.333
$ws.Cells.Item(98, 12).Value = 15999.5  # L98: was 15766.333
$ws.Cells.Item(98, 14).Value = -21989.5  # N98: was -21756.333

# Row 107 (sheet WVR)
$ws.Cells.Item(107, 8).Value = 669.11536  # H107: was 717.7143
$ws.Cells.Item(107, 9).Value = 716.26666  # I107: was 696.4375
$ws.Cells.Item(107, 10).Value = 604.8182  # J107: was 785.8
$ws.Cells.Item(107, 11).Value = 2148.79998  # K107: was 2089.3125
$ws.Cells.Item(107, 12).Value = 1814.4546  # L107: was 2357.4
$ws.Cells.Item(107, 13).Value = -228.7999799999998  # M107: was -169.3125
$ws.Cells.Item(107, 14).Value = -5654.4546  # N107: was -6197.4

# Row 122 (sheet WVR)
$ws.Cells.Item(122, 8).Value = 2285.4167  # H122: was 1641.3889
$ws.Cells.Item(122, 9).Value = 2038.6364  # I122: was 1443.8235
$ws.Cells.Item(122, 11).Value = 6115.9092  # K122: was 4331.470499999999
$ws.Cells.Item(122, 13).Value = -3665.9092  # M122: was -1881.470499999999
